# B6-PowerPoint.pptx edit:
#  1. Three tables (slides 14, 15, 16) switch from the custom "Table_0"
#     style to the built-in "No Style, Table Grid" style.
#  2. The presentation's applied colour theme changes from the
#     "Integral" / Red Violet scheme to the standard Office colour
#     scheme (Design > Variants > Colors > "Office").

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables -----------------------------------------
$noStyleTableGrid = "{AC4A7E81-45A9-4FD1-BBA0-68BF6FC0D2DC}"
foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    $tableShape = $slide.Shapes.Item(1)
    if ($tableShape.HasTable) {
        $tableShape.Table.ApplyStyle($noStyleTableGrid)
    }
}

# --- 2. Swap the colour scheme from "Red Violet" to "Office" --------------
# PpThemeColorSchemeIndex order: Dark1, Light1, Dark2, Light2,
# Accent1-6, Hyperlink, FollowedHyperlink.
$officeColors = @(
    0,          # Dark1   000000
    16777215,   # Light1  FFFFFF
    6968388,    # Dark2   44546A
    15132391,   # Light2  E7E6E6
    13998939,   # Accent1 5B9BD5
    3243501,    # Accent2 ED7D31
    10855845,   # Accent3 A5A5A5
    49407,      # Accent4 FFC000
    12874308,   # Accent5 4472C4
    4697456,    # Accent6 70AD47
    12673797,   # Hyperlink 0563C1
    7491477     # FollowedHyperlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
